$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# A new formula entry ("UNION_INDICATORS_PT4") was added to the library.
# Insert a fresh row at position 8 -- this pushes the existing rows 8..102
# down to 9..103 -- and fill it in with the new entry's data.
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "CREATE/MODIFY"
$ws.Range("B8").Value = "LIB_EWS_RETAIL"
$ws.Range("C8").Value = "UNION_INDICATORS_PT4"
$ws.Range("E8").Value = "String"
$ws.Range("F8").Value = "String"

# Reflect the edit in the sheet view: the cursor ends up on the newly
# inserted cell instead of the previous scroll position further down.
$ws.Activate()
$ws.Range("C8").Select()
